$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: deaths_demo
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("deaths_demo")

$ws1.Range("B2").Value = 154
$ws1.Range("C2").Value = 6

$ws1.Range("B3").Value = 1933
$ws1.Range("C3").Value = 134
$ws1.Range("D3").Value = 31
$ws1.Range("E3").Value = 3
$ws1.Range("F3").Value = 0.2

$ws1.Range("B4").Value = 836
$ws1.Range("C4").Value = 195
$ws1.Range("D4").Value = 66
$ws1.Range("E4").Value = 14
$ws1.Range("F4").Value = 1.7

$ws1.Range("B5").Value = 453
$ws1.Range("C5").Value = 176
$ws1.Range("D5").Value = 68
$ws1.Range("E5").Value = 43
$ws1.Range("F5").Value = 9.5

$ws1.Range("B6").Value = 374
$ws1.Range("C6").Value = 155
$ws1.Range("D6").Value = 55
$ws1.Range("E6").Value = 80
$ws1.Range("F6").Value = 21.4

$ws1.Range("B7").Value = 436
$ws1.Range("C7").Value = 121
$ws1.Range("D7").Value = 17
$ws1.Range("E7").Value = 131
$ws1.Range("F7").Value = 30

# ---------------------------------------------------------------------------
# Sheet 2: ethnicities
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ethnicities")

$ws2.Range("B2").Value = 1961
$ws2.Range("C2").Value = 1397
$ws2.Range("D2").Value = 429
$ws2.Range("E2").Value = 105
$ws2.Range("F2").Value = 116
$ws2.Range("G2").Value = 178

$ws2.Range("B3").Value = 403
$ws2.Range("C3").Value = 259
$ws2.Range("D3").Value = 80
$ws2.Range("E3").Value = 14
$ws2.Range("F3").Value = 25
$ws2.Range("G3").Value = 6

$ws2.Range("B4").Value = 109
$ws2.Range("C4").Value = 95
$ws2.Range("D4").Value = 24
$ws2.Range("E4").Value = 4
$ws2.Range("F4").Value = 4
$ws2.Range("G4").Value = 1

$ws2.Range("B5").Value = 200
$ws2.Range("C5").Value = 42
$ws2.Range("D5").Value = 17
$ws2.Range("E5").Value = 4
$ws2.Range("F5").Value = 3
$ws2.Range("G5").Value = 5

$ws2.Range("B6").Value = 91
$ws2.Range("C6").Value = 26
$ws2.Range("D6").Value = 13
$ws2.Range("E6").Value = 1
$ws2.Range("F6").Value = 2
$ws2.Range("G6").Value = 3

# Row 7 ("Out-of-Hospital Deaths") is a shared formula "=Bx5-Bx6" that
# recalculates automatically now that its precedents (rows 5 and 6) changed;
# no direct edit to the formula cells themselves is required or desired.

# ---------------------------------------------------------------------------
# Sheet 3: prop
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("prop")

$ws3.Range("B2").Value = 46.846631629240328
$ws3.Range("C2").Value = 33.373148590539898
$ws3.Range("D2").Value = 10.248447204968944
$ws3.Range("E2").Value = 2.508361204013378
$ws3.Range("F2").Value = 2.7711419015766841
$ws3.Range("G2").Value = 4.2522694696607743

$ws3.Range("B3").Value = 51.207115628970776
$ws3.Range("C3").Value = 32.909783989834821
$ws3.Range("D3").Value = 10.165184243964422
$ws3.Range("E3").Value = 1.7789072426937738
$ws3.Range("F3").Value = 3.1766200762388821
$ws3.Range("G3").Value = 0.76238881829733163

$ws3.Range("B4").Value = 45.991561181434598
$ws3.Range("C4").Value = 40.084388185654007
$ws3.Range("D4").Value = 10.126582278481013
$ws3.Range("E4").Value = 1.6877637130801686
$ws3.Range("F4").Value = 1.6877637130801686
$ws3.Range("G4").Value = 0.42194092827004215

$ws3.Range("B5").Value = 73.800738007380076
$ws3.Range("C5").Value = 15.498154981549817
$ws3.Range("D5").Value = 6.2730627306273057
$ws3.Range("E5").Value = 1.4760147601476015
$ws3.Range("F5").Value = 1.107011070110701
$ws3.Range("G5").Value = 1.8450184501845017

$ws3.Range("B6").Value = 66.911764705882348
$ws3.Range("C6").Value = 19.117647058823529
$ws3.Range("D6").Value = 9.5588235294117645
$ws3.Range("E6").Value = 0.73529411764705876
$ws3.Range("F6").Value = 1.4705882352941175
$ws3.Range("G6").Value = 2.2058823529411766

$ws3.Range("B7").Value = 80.740740740740748
$ws3.Range("C7").Value = 11.851851851851853
$ws3.Range("D7").Value = 2.9629629629629632
$ws3.Range("E7").Value = 2.2222222222222223
$ws3.Range("F7").Value = 0.74074074074074081
$ws3.Range("G7").Value = 1.4814814814814816

# New empty, bold-formatted cells H12:H17 (no values, formatting only)
$ws3.Range("H12:H17").Font.Bold = $true

$ws3.Activate()
[void]$ws3.Range("G21").Select()
